$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText, $matchWholeWord) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $matchWholeWord, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $findText"
        return $false
    }
    $r.Text = $replaceText
    return $true
}

$apos = [char]0x2019

# 1. "Solution – Take the Bird come back..." -> add commas, drop proofErr around "come"
$ok = ReplaceText "Solution – Take the Bird come back and get the cat and then the seed " "Solution – Take the Bird, come back and get the cat, and then the seed " $false
Write-Host "1:" $ok

# 2. "Whats the min a amount" -> "What's the min  amount" (merge runs, drop 3 proofErr pairs)
$ok = ReplaceText "a) Whats the min a amount of socks you can grab to get a match " ("a) What" + $apos + "s the min  amount of socks you can grab to get a match ") $false
Write-Host "2:" $ok

# 3. "b) are they already mated up" -> "b) Are ..." (include "b) " to drop both proofErr tags)
$ok = ReplaceText "b) are they already mated up " "b) Are they already mated up " $false
Write-Host "3:" $ok

# 4. "a) grab four socks..." -> "a) Grab four socks..."
$ok = ReplaceText "a) grab four socks to get mates and 12 socks to get on of each color " "a) Grab four socks to get mates and 12 socks to get on of each color " $false
Write-Host "4:" $ok

# 5. "b) if you grab too few..." -> "b) If you grab too few..."
$ok = ReplaceText "b) if you grab too few you wont get one of each or a match " "b) If you grab too few you wont get one of each or a match " $false
Write-Host "5:" $ok

# 6. "Which finger will she stop on " -> "...on? " (include "a) " context to drop both proofErr tags)
$ok = ReplaceText "a) Which finger will she stop on " "a) Which finger will she stop on? " $false
Write-Host "6:" $ok

# 7. "She stops on a differnet finger" -> "...different..."
$ok = ReplaceText "She stops on a differnet finger each time " "She stops on a different finger each time " $false
Write-Host "7:" $ok

# 8. "count out on hand" -> "Count out on hand" ; can't safely cross the <w:tab/> boundary,
#    so only the gramEnd (trailing marker) is safely absorbed; gramStart remains orphaned.
$ok = ReplaceText "count out on hand" "Count out on hand" $false
Write-Host "8:" $ok

# 9. "4)  Evaluate each potential Solution " -> "4) Evaluate each potential Solution " (fix double space, drop proofErr)
$ok = ReplaceText "4)  Evaluate each potential Solution " "4) Evaluate each potential Solution " $false
Write-Host "9:" $ok

# 10. "5)  Choose a Solution and Develop a Plan to Implement it " -> single space, drop proofErr
$ok = ReplaceText "5)  Choose a Solution and Develop a Plan to Implement it " "5) Choose a Solution and Develop a Plan to Implement it " $false
Write-Host "10:" $ok

Write-Host "Done so far"
